{"js": "const pairs = [\n  [\"14+47=61\", \"98-64=34\"],\n  [\"93-5=88\", \"67-56=11\"],\n  [\"67-19=48\", \"27-26=1\"],\n  [\"11-0=11\", \"57-20=37\"],\n  [\"75-21=54\", \"9+88=97\"],\n  [\"74-60=14\", \"24+4=28\"],\n  [\"95-89=6\", \"3+41=44\"],\n  [\"97-24=73\", \"65+6=71\"],\n  [\"34-28=6\", \"68+4=72\"],\n  [\"16+81=97\", \"89-57=32\"],\n  [\"1+47=48\", \"11+0=11\"],\n  [\"75-42=33\", \"90-24=66\"],\n  [\"33-5=28\", \"54-37=17\"],\n  [\"60-35=25\", \"83-9=74\"],\n  [\"53+28=81\", \"52+8=60\"],\n  [\"24-11=13\", \"70+23=93\"],\n  [\"31+0=31\", \"68+2=70\"],\n  [\"29-14=15\", \"52-51=1\"],\n  [\"44-44=0\", \"92-46=46\"],\n  [\"67+10=77\", \"92+2=94\"],\n  [\"79+5=84\", \"99-23=76\"],\n  [\"11+33=44\", \"93-20=73\"],\n  [\"14+4=18\", \"90-88=2\"],\n  [\"33+42=75\", \"33+30=63\"],\n  [\"96-29=67\", \"76-30=46\"],\n  [\"3+94=97\", \"71-52=19\"],\n  [\"71+15=86\", \"23+63=86\"],\n  [\"51+10=61\", \"10+25=35\"],\n  [\"13+71=84\", \"61-22=39\"],\n  [\"61+23=84\", \"84-23=61\"],\n  [\"75+8=83\", \"92-37=55\"],\n  [\"1+33=34\", \"65-43=22\"],\n  [\"87-73=14\", \"20+36=56\"],\n  [\"87-59=28\", \"77-8=69\"],\n  [\"7-0=7\", \"75-14=61\"],\n  [\"46-22=24\", \"42-14=28\"],\n  [\"48+35=83\", \"49+22=71\"],\n  [\"97-32=65\", \"76-63=13\"],\n  [\"19+24=43\", \"44-36=8\"],\n  [\"97-63=34\", \"42+12=54\"],\n  [\"9+59=68\", \"55-19=36\"],\n  [\"37+30=67\", \"43-2=41\"],\n  [\"76-73=3\", \"51+36=87\"],\n  [\"42-13=29\", \"47-13=34\"],\n  [\"11+80=91\", \"49+21=70\"],\n  [\"18+5=23\", \"75-70=5\"],\n  [\"9+40=49\", \"52+18=70\"],\n  [\"24+20=44\", \"67-32=35\"],\n  [\"47+33=80\", \"46-5=41\"],\n  [\"14+15=29\", \"75+11=86\"],\n  [\"21+21=42\", \"59-21=38\"],\n  [\"75-56=19\", \"30+62=92\"],\n  [\"98-1=97\", \"34+37=71\"],\n  [\"92-18=74\", \"46+37=83\"],\n  [\"77-45=32\", \"19+63=82\"],\n  [\"4+83=87\", \"35-27=8\"],\n  [\"16+9=25\", \"7+83=90\"],\n  [\"63+9=72\", \"61-2=59\"],\n  [\"83-74=9\", \"9+71=80\"],\n  [\"54+20=74\", \"75-70=5\"],\n  [\"91-89=2\", \"30-21=9\"],\n  [\"95-43=52\", \"98-44=54\"],\n  [\"18-14=4\", \"1+6=7\"],\n  [\"93-24=69\", \"66-63=3\"],\n  [\"99-17=82\", \"81-58=23\"],\n  [\"40-24=16\", \"76-47=29\"],\n  [\"23+72=95\", \"13+32=45\"],\n  [\"36+9=45\", \"6+75=81\"],\n  [\"74+14=88\", \"80-34=46\"],\n  [\"86-27=59\", \"84-53=31\"],\n  [\"11+18=29\", \"45+20=65\"],\n  [\"37+17=54\", \"62-50=12\"],\n  [\"98-27=71\", \"14+81=95\"],\n  [\"43+23=66\", \"61-32=29\"],\n  [\"15+21=36\", \"59+26=85\"],\n  [\"47-9=38\", \"81+9=90\"],\n  [\"28+56=84\", \"10+55=65\"],\n  [\"77-52=25\", \"20-16=4\"],\n  [\"81+0=81\", \"13+76=89\"],\n  [\"27+65=92\", \"49-16=33\"],\n  [\"92-40=52\", \"49+22=71\"],\n  [\"68-51=17\", \"64-53=11\"],\n  [\"56+0=56\", \"9+79=88\"],\n  [\"19+67=86\", \"62+1=63\"],\n  [\"98-6=92\", \"12+51=63\"],\n  [\"89-64=25\", \"36-15=21\"],\n  [\"44+26=70\", \"81-21=60\"],\n  [\"83-62=21\", \"94-63=31\"],\n  [\"58-56=2\", \"23+16=39\"],\n  [\"97-14=83\", \"53-4=49\"],\n  [\"33+18=51\", \"26-19=7\"],\n  [\"11+17=28\", \"64+29=93\"],\n  [\"67-9=58\", \"13+63=76\"],\n  [\"24+10=34\", \"25+46=71\"],\n  [\"35+63=98\", \"82-50=32\"],\n  [\"62+11=73\", \"92-88=4\"],\n  [\"46+11=57\", \"95-57=38\"],\n  [\"69+8=77\", \"53-52=1\"],\n  [\"85-70=15\", \"56+37=93\"],\n  [\"0+91=91\", \"97-73=24\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No tables found in document body\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst columnCount = Math.round(pairs.length / table.rowCount);\n\n// Collect the paragraph range for every cell, in row-major (document) order.\nconst cellRanges = [];\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < columnCount; c++) {\n    const cell = table.getCell(r, c);\n    const para = cell.body.paragraphs.getFirst();\n    cellRanges.push(para.getRange());\n  }\n}\n\nif (cellRanges.length !== pairs.length) {\n  throw new Error(\n    \"Cell count (\" + cellRanges.length + \") does not match expected pair count (\" + pairs.length + \")\"\n  );\n}\n\nfor (let i = 0; i < cellRanges.length; i++) {\n  cellRanges[i].load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < cellRanges.length; i++) {\n  const expected = pairs[i][0];\n  const replacement = pairs[i][1];\n  const actual = cellRanges[i].text;\n  if (actual !== expected) {\n    throw new Error(\n      \"Cell \" + i + \" text mismatch: expected '\" + expected + \"' but found '\" + actual + \"'\"\n    );\n  }\n  cellRanges[i].insertText(replacement, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$pairs = @(\n  @(\"14+47=61\", \"98-64=34\"),\n  @(\"93-5=88\", \"67-56=11\"),\n  @(\"67-19=48\", \"27-26=1\"),\n  @(\"11-0=11\", \"57-20=37\"),\n  @(\"75-21=54\", \"9+88=97\"),\n  @(\"74-60=14\", \"24+4=28\"),\n  @(\"95-89=6\", \"3+41=44\"),\n  @(\"97-24=73\", \"65+6=71\"),\n  @(\"34-28=6\", \"68+4=72\"),\n  @(\"16+81=97\", \"89-57=32\"),\n  @(\"1+47=48\", \"11+0=11\"),\n  @(\"75-42=33\", \"90-24=66\"),\n  @(\"33-5=28\", \"54-37=17\"),\n  @(\"60-35=25\", \"83-9=74\"),\n  @(\"53+28=81\", \"52+8=60\"),\n  @(\"24-11=13\", \"70+23=93\"),\n  @(\"31+0=31\", \"68+2=70\"),\n  @(\"29-14=15\", \"52-51=1\"),\n  @(\"44-44=0\", \"92-46=46\"),\n  @(\"67+10=77\", \"92+2=94\"),\n  @(\"79+5=84\", \"99-23=76\"),\n  @(\"11+33=44\", \"93-20=73\"),\n  @(\"14+4=18\", \"90-88=2\"),\n  @(\"33+42=75\", \"33+30=63\"),\n  @(\"96-29=67\", \"76-30=46\"),\n  @(\"3+94=97\", \"71-52=19\"),\n  @(\"71+15=86\", \"23+63=86\"),\n  @(\"51+10=61\", \"10+25=35\"),\n  @(\"13+71=84\", \"61-22=39\"),\n  @(\"61+23=84\", \"84-23=61\"),\n  @(\"75+8=83\", \"92-37=55\"),\n  @(\"1+33=34\", \"65-43=22\"),\n  @(\"87-73=14\", \"20+36=56\"),\n  @(\"87-59=28\", \"77-8=69\"),\n  @(\"7-0=7\", \"75-14=61\"),\n  @(\"46-22=24\", \"42-14=28\"),\n  @(\"48+35=83\", \"49+22=71\"),\n  @(\"97-32=65\", \"76-63=13\"),\n  @(\"19+24=43\", \"44-36=8\"),\n  @(\"97-63=34\", \"42+12=54\"),\n  @(\"9+59=68\", \"55-19=36\"),\n  @(\"37+30=67\", \"43-2=41\"),\n  @(\"76-73=3\", \"51+36=87\"),\n  @(\"42-13=29\", \"47-13=34\"),\n  @(\"11+80=91\", \"49+21=70\"),\n  @(\"18+5=23\", \"75-70=5\"),\n  @(\"9+40=49\", \"52+18=70\"),\n  @(\"24+20=44\", \"67-32=35\"),\n  @(\"47+33=80\", \"46-5=41\"),\n  @(\"14+15=29\", \"75+11=86\"),\n  @(\"21+21=42\", \"59-21=38\"),\n  @(\"75-56=19\", \"30+62=92\"),\n  @(\"98-1=97\", \"34+37=71\"),\n  @(\"92-18=74\", \"46+37=83\"),\n  @(\"77-45=32\", \"19+63=82\"),\n  @(\"4+83=87\", \"35-27=8\"),\n  @(\"16+9=25\", \"7+83=90\"),\n  @(\"63+9=72\", \"61-2=59\"),\n  @(\"83-74=9\", \"9+71=80\"),\n  @(\"54+20=74\", \"75-70=5\"),\n  @(\"91-89=2\", \"30-21=9\"),\n  @(\"95-43=52\", \"98-44=54\"),\n  @(\"18-14=4\", \"1+6=7\"),\n  @(\"93-24=69\", \"66-63=3\"),\n  @(\"99-17=82\", \"81-58=23\"),\n  @(\"40-24=16\", \"76-47=29\"),\n  @(\"23+72=95\", \"13+32=45\"),\n  @(\"36+9=45\", \"6+75=81\"),\n  @(\"74+14=88\", \"80-34=46\"),\n  @(\"86-27=59\", \"84-53=31\"),\n  @(\"11+18=29\", \"45+20=65\"),\n  @(\"37+17=54\", \"62-50=12\"),\n  @(\"98-27=71\", \"14+81=95\"),\n  @(\"43+23=66\", \"61-32=29\"),\n  @(\"15+21=36\", \"59+26=85\"),\n  @(\"47-9=38\", \"81+9=90\"),\n  @(\"28+56=84\", \"10+55=65\"),\n  @(\"77-52=25\", \"20-16=4\"),\n  @(\"81+0=81\", \"13+76=89\"),\n  @(\"27+65=92\", \"49-16=33\"),\n  @(\"92-40=52\", \"49+22=71\"),\n  @(\"68-51=17\", \"64-53=11\"),\n  @(\"56+0=56\", \"9+79=88\"),\n  @(\"19+67=86\", \"62+1=63\"),\n  @(\"98-6=92\", \"12+51=63\"),\n  @(\"89-64=25\", \"36-15=21\"),\n  @(\"44+26=70\", \"81-21=60\"),\n  @(\"83-62=21\", \"94-63=31\"),\n  @(\"58-56=2\", \"23+16=39\"),\n  @(\"97-14=83\", \"53-4=49\"),\n  @(\"33+18=51\", \"26-19=7\"),\n  @(\"11+17=28\", \"64+29=93\"),\n  @(\"67-9=58\", \"13+63=76\"),\n  @(\"24+10=34\", \"25+46=71\"),\n  @(\"35+63=98\", \"82-50=32\"),\n  @(\"62+11=73\", \"92-88=4\"),\n  @(\"46+11=57\", \"95-57=38\"),\n  @(\"69+8=77\", \"53-52=1\"),\n  @(\"85-70=15\", \"56+37=93\"),\n  @(\"0+91=91\", \"97-73=24\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nif ($rowCount * $colCount -ne $pairs.Count) {\n    throw \"Cell count ($($rowCount * $colCount)) does not match expected pair count ($($pairs.Count))\"\n}\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $pair = $pairs[$idx]\n        $expected = $pair[0]\n        $replacement = $pair[1]\n\n        $cell = $t.Cell($r, $c)\n        $cellRange = $cell.Range\n        $cellRange.MoveEnd(1, -1) | Out-Null\n\n        $actual = $cellRange.Text\n        if ($actual -ne $expected) {\n            throw \"Cell ($r,$c) text mismatch: expected '$expected' but found '$actual'\"\n        }\n\n        $cellRange.Text = $replacement\n\n        $idx = $idx + 1\n    }\n}\n"}
